$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1378.2759
$ws.Range("J17").Value = 1406.2963
$ws.Range("L17").Value = 4218.8889
$ws.Range("N17").Value = -4554.8889

$ws.Range("H33").Value = 517.3333
$ws.Range("I33").Value = 532.3158
$ws.Range("K33").Value = 532.3158
$ws.Range("M33").Value = -303.3158

$ws.Range("H116").Value = 1600
$ws.Range("I116").Value = 1600
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 1600
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = 1842
$ws.Range("N116").Value = ""

$ws.Range("H125").Value = 115389500
$ws.Range("I125").Value = 90912690
$ws.Range("K125").Value = 818214210
$ws.Range("M125").Value = -818211750

$ws.Range("H138").Value = 2626.3333
$ws.Range("I138").Value = 2321
$ws.Range("K138").Value = 6963
$ws.Range("M138").Value = -1823

$ws.Range("H141").Value = 1833.5
$ws.Range("I141").Value = 1833.5
$ws.Range("K141").Value = 5500.5
$ws.Range("M141").Value = -320.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H17").Value = 17500
$ws.Range("I17").Value = 15000
$ws.Range("J17").Value = 20000
$ws.Range("K17").Value = 15000
$ws.Range("L17").Value = 20000
$ws.Range("M17").Value = -14827
$ws.Range("N17").Value = -20346

$ws.Range("H32").Value = 2447.5293
$ws.Range("I32").Value = 2513.4243
$ws.Range("K32").Value = 2513.4243
$ws.Range("M32").Value = -2226.4243

$ws.Range("H63").Value = 17179.2
$ws.Range("I63").Value = 18248.666
$ws.Range("K63").Value = 18248.666
$ws.Range("M63").Value = -17562.666

$ws.Range("H66").Value = 17179.2
$ws.Range("I66").Value = 18248.666
$ws.Range("K66").Value = 91243.33
$ws.Range("M66").Value = -87811.33

$ws.Range("H88").Value = 686.8333
$ws.Range("I88").Value = 708
$ws.Range("J88").Value = 665.6667
$ws.Range("K88").Value = 708
$ws.Range("L88").Value = 665.6667
$ws.Range("M88").Value = -302
$ws.Range("N88").Value = -1477.6667

$ws.Range("H91").Value = 686.8333
$ws.Range("I91").Value = 708
$ws.Range("J91").Value = 665.6667
$ws.Range("K91").Value = 708
$ws.Range("L91").Value = 665.6667
$ws.Range("M91").Value = 696
$ws.Range("N91").Value = -3473.6667

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H88").Value = 14399.75
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 14399.75
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 14399.75
$ws.Range("M88").Value = ""
$ws.Range("N88").Value = -15211.75

$ws.Range("H91").Value = 14399.75
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 14399.75
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 14399.75
$ws.Range("M91").Value = ""
$ws.Range("N91").Value = -17207.75

$ws.Range("H113").Value = 14850
$ws.Range("I113").Value = 14850
$ws.Range("K113").Value = 14850
$ws.Range("M113").Value = -12680

$ws.Range("H122").Value = 70780
$ws.Range("J122").Value = 70780
$ws.Range("L122").Value = 70780
$ws.Range("N122").Value = -80580

$ws.Range("H134").Value = 3773.1875
$ws.Range("I134").Value = 3619.4285
$ws.Range("K134").Value = 10858.2855
$ws.Range("M134").Value = -8323.2855

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3988.5557
$ws.Range("I16").Value = 3419.8
$ws.Range("J16").Value = 4699.5
$ws.Range("K16").Value = 3419.8
$ws.Range("L16").Value = 4699.5
$ws.Range("M16").Value = -3132.8
$ws.Range("N16").Value = -5273.5

$ws.Range("H31").Value = 3153.3333
$ws.Range("J31").Value = 3497
$ws.Range("L31").Value = 3497
$ws.Range("N31").Value = -4087

$ws.Range("H34").Value = 3153.3333
$ws.Range("J34").Value = 3497
$ws.Range("L34").Value = 3497
$ws.Range("N34").Value = -3901

$ws.Range("H81").Value = 99999
$ws.Range("J81").Value = 99999
$ws.Range("L81").Value = 99999
$ws.Range("N81").Value = -101995

$ws.Range("H84").Value = 99999
$ws.Range("J84").Value = 99999
$ws.Range("L84").Value = 299997
$ws.Range("N84").Value = -309981

$ws.Range("H113").Value = 3988.5557
$ws.Range("I113").Value = 3419.8
$ws.Range("J113").Value = 4699.5
$ws.Range("K113").Value = 3419.8
$ws.Range("L113").Value = 4699.5
$ws.Range("M113").Value = -1249.8
$ws.Range("N113").Value = -9039.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 15599.833
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 15599.833
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 46799.499
$ws.Range("M70").Value = ""
$ws.Range("N70").Value = -47429.499

$ws.Range("H73").Value = 15599.833
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 15599.833
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 46799.499
$ws.Range("M73").Value = ""
$ws.Range("N73").Value = -48983.499

$ws.Range("H107").Value = 966.6667
$ws.Range("J107").Value = 950
$ws.Range("L107").Value = 2850
$ws.Range("N107").Value = -6690

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 640781
$ws.Range("I14").Value = 678438.7
$ws.Range("J14").Value = 600
$ws.Range("K14").Value = 678438.7
$ws.Range("L14").Value = 600
$ws.Range("M14").Value = -678270.7
$ws.Range("N14").Value = -936

$ws.Range("H46").Value = 1000
$ws.Range("I46").Value = 1000
$ws.Range("K46").Value = 1000
$ws.Range("M46").Value = -844

$ws.Range("H80").Value = 2116.3333
$ws.Range("I80").Value = 2116.3333
$ws.Range("K80").Value = 2116.3333
$ws.Range("M80").Value = -1118.3333

$ws.Range("H83").Value = 2116.3333
$ws.Range("I83").Value = 2116.3333
$ws.Range("K83").Value = 10581.6665
$ws.Range("M83").Value = -5589.666499999999

$ws.Range("H98").Value = 27989.572
$ws.Range("J98").Value = 27989.572
$ws.Range("L98").Value = 27989.572
$ws.Range("N98").Value = -33979.572

$ws.Range("H107").Value = 2173.2778
$ws.Range("I107").Value = 1419.7
$ws.Range("J107").Value = 3115.25
$ws.Range("K107").Value = 1419.7
$ws.Range("L107").Value = 3115.25
$ws.Range("M107").Value = 500.3
$ws.Range("N107").Value = -6955.25

$ws.Range("H113").Value = 1399
$ws.Range("I113").Value = 865.3333
$ws.Range("J113").Value = 3000
$ws.Range("K113").Value = 865.3333
$ws.Range("L113").Value = 3000
$ws.Range("M113").Value = 1304.6667
$ws.Range("N113").Value = -7340

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5499.28
$ws.Range("I7").Value = 2049.5
$ws.Range("J7").Value = 7799.1333
$ws.Range("K7").Value = 2049.5
$ws.Range("L7").Value = 7799.1333
$ws.Range("M7").Value = -1937.5
$ws.Range("N7").Value = -8023.1333

$ws.Range("H35").Value = 643.3333
$ws.Range("I35").Value = 643.3333
$ws.Range("K35").Value = 643.3333
$ws.Range("M35").Value = -307.3333

$ws.Range("H40").Value = 3139.2
$ws.Range("I40").Value = 2465.75
$ws.Range("J40").Value = 5833
$ws.Range("K40").Value = 2465.75
$ws.Range("L40").Value = 5833
$ws.Range("M40").Value = -2329.75
$ws.Range("N40").Value = -6105

$ws.Range("H61").Value = 4094.6
$ws.Range("J61").Value = 4333.3335
$ws.Range("L61").Value = 4333.3335
$ws.Range("N61").Value = -4737.3335

$ws.Range("H82").Value = 1315
$ws.Range("I82").Value = 1571
$ws.Range("J82").Value = 867
$ws.Range("K82").Value = 1571
$ws.Range("L82").Value = 867
$ws.Range("M82").Value = -1210
$ws.Range("N82").Value = -1589

$ws.Range("H85").Value = 1315
$ws.Range("I85").Value = 1571
$ws.Range("J85").Value = 867
$ws.Range("K85").Value = 1571
$ws.Range("L85").Value = 867
$ws.Range("M85").Value = -323
$ws.Range("N85").Value = -3363

$ws.Range("H113").Value = 4094.6
$ws.Range("J113").Value = 4333.3335
$ws.Range("L113").Value = 4333.3335
$ws.Range("N113").Value = -8673.333500000001

$ws.Range("H122").Value = 5972.6333
$ws.Range("I122").Value = 5791.28
$ws.Range("K122").Value = 17373.84
$ws.Range("M122").Value = -14923.84

$ws.Range("H126").Value = 5499.28
$ws.Range("I126").Value = 2049.5
$ws.Range("J126").Value = 7799.1333
$ws.Range("K126").Value = 6148.5
$ws.Range("L126").Value = 23397.3999
$ws.Range("M126").Value = -3678.5
$ws.Range("N126").Value = -28337.3999

$ws.Range("H136").Value = 2120
$ws.Range("I136").Value = 2120
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 6360
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -3810
$ws.Range("N136").Value = ""

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 89064
$ws.Range("I75").Value = 89018
$ws.Range("K75").Value = 89018
$ws.Range("M75").Value = -88082

$ws.Range("H78").Value = 89064
$ws.Range("I78").Value = 89018
$ws.Range("K78").Value = 267054
$ws.Range("M78").Value = -262374

$ws.Range("H107").Value = 1495
$ws.Range("I107").Value = 1062.5714
$ws.Range("K107").Value = 3187.7142
$ws.Range("M107").Value = -1267.7142
